$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kayitlar")

# Insert a new row above the existing data row (row 2), shifting it down to row 3
$ws.Rows.Item(2).Insert()

# Keep the new row's cells as plain text, matching the rest of the sheet
$ws.Range("A2:G2").NumberFormat = "@"

# New record (Kayıt No 3, tarih 2025-07-15) goes into the now-empty row 2
$ws.Range("A2").Value = "3"
$ws.Range("B2").Value = "2025-07-15"
$ws.Range("C2").Value = "İlçe"
$ws.Range("D2").Value = "2"
$ws.Range("E2").Value = "2"
$ws.Range("F2").Value = "Cins D."
$ws.Range("G2").Value = "Gökhan ELGÜL"
